# Tilaukset logataan myös tekstifiluun
#
# The order/address block (name, street, city) that used to be a run of
# plain top-level paragraphs is moved into the first cell of a new 2x1
# table, font-sized at 12pt (sz 24 half-points), and a second, empty
# cell is added next to it. The trailing "undefined undefined" filler
# paragraph is dropped. A blank paragraph carrying the section
# properties is left at the top of the body (this is where the old
# section break now lives now that the content moved into the table).
#
# We rebuild the whole main story (everything except the document's
# final sectPr, which Range.InsertXML leaves untouched when the target
# range doesn't include it) from a literal OOXML fragment via
# Range.InsertXML — this gives byte-exact control over the resulting
# markup instead of approximating it through many small object-model
# calls.

$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:sectPr>
                <w:pgSz w:w="11906" w:h="16838" w:orient="portrait"/>
                <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="708" w:footer="708" w:gutter="0"/>
                <w:pgNumType/>
                <w:docGrid w:linePitch="360"/>
              </w:sectPr>
            </w:pPr>
          </w:p>
          <w:tbl>
            <w:tblPr>
              <w:tblW w:type="auto" w:w="100"/>
              <w:tblBorders>
                <w:top w:val="single" w:color="auto" w:sz="4"/>
                <w:left w:val="single" w:color="auto" w:sz="4"/>
                <w:bottom w:val="single" w:color="auto" w:sz="4"/>
                <w:right w:val="single" w:color="auto" w:sz="4"/>
                <w:insideH w:val="single" w:color="auto" w:sz="4"/>
                <w:insideV w:val="single" w:color="auto" w:sz="4"/>
              </w:tblBorders>
            </w:tblPr>
            <w:tblGrid>
              <w:gridCol w:w="100"/>
              <w:gridCol w:w="100"/>
            </w:tblGrid>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:type="pct" w:w="50%"/>
                  <w:tcMar>
                    <w:top w:type="dxa" w:w="100"/>
                    <w:bottom w:type="dxa" w:w="100"/>
                  </w:tcMar>
                </w:tcPr>
                <w:p>
                  <w:r>
                    <w:rPr>
                      <w:b/>
                      <w:bCs/>
                      <w:sz w:val="24"/>
                      <w:szCs w:val="24"/>
                    </w:rPr>
                    <w:t xml:space="preserve">John Doe</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:r>
                    <w:rPr>
                      <w:sz w:val="24"/>
                      <w:szCs w:val="24"/>
                    </w:rPr>
                    <w:t xml:space="preserve">123 Main St</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:r>
                    <w:rPr>
                      <w:sz w:val="24"/>
                      <w:szCs w:val="24"/>
                    </w:rPr>
                    <w:t xml:space="preserve">Anytown 12345</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:pPr>
                    <w:spacing w:after="200"/>
                  </w:pPr>
                </w:p>
              </w:tc>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:type="pct" w:w="50%"/>
                  <w:tcMar>
                    <w:top w:type="dxa" w:w="100"/>
                    <w:bottom w:type="dxa" w:w="100"/>
                  </w:tcMar>
                </w:tcPr>
                <w:p/>
              </w:tc>
            </w:tr>
          </w:tbl>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($xml)
